$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.624.91'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.221.11'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.52'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.06'
$ws.Range('E6').Value = '  -5.93%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.565'
$ws.Range('E7').Value = '  -3.73%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  -7.49%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.56'
$ws.Range('E10').Value = '  -8.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  -3.14%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.32'
$ws.Range('E12').Value = '  -7.30%  '
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.558.83'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.220.63'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.832'
$ws.Range('E16').Value = '  -5.08%  '
$ws.Range('E17').Value = '  -4.23%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.498.51'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.95'
$ws.Range('E19').Value = '  -10.42%  '
$ws.Range('E20').Value = '  -4.60%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.26'
$ws.Range('E21').Value = '  -6.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '64.90'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '234.75'
$ws.Range('E23').Value = '  -1.50%  '
$ws.Range('E24').Value = '  -8.39%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  -8.52%  '
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.90'
$ws.Range('E27').Value = '  -3.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  -1.37%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '36.17'
$ws.Range('E29').Value = '  -8.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.92'
$ws.Range('E30').Value = '  -9.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '157.69'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.75'
$ws.Range('E32').Value = '  -3.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0820'
$ws.Range('E33').Value = '  -6.91%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.65'
$ws.Range('E34').Value = '  -2.68%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.08'
$ws.Range('E35').Value = '  -6.95%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('E37').Value = '  -10.06%  '
$ws.Range('E38').Value = '  -3.74%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '15.41'
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.50'
$ws.Range('E40').Value = '  -9.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.99'
$ws.Range('E41').Value = '  -12.18%  '
$ws.Range('E42').Value = '  -7.24%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.691.32'
$ws.Range('E44').Value = '  -4.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '83.98'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('E46').Value = '  -8.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.09'
$ws.Range('E47').Value = '  -6.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '100.81'
$ws.Range('E48').Value = '  -3.77%  '
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '70.49'
$ws.Range('E50').Value = '  -5.80%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '55.50'
$ws.Range('E51').Value = '  -6.90%  '
